# Project Plan & Technical debt.xlsx - "Add files via upload" edit
#
# Summary of the change (reconstructed from the OOXML diff):
#  1. The "Activity-tasks map list" sheet is renamed to a name that is
#     literally wrapped in double quotes: "Activity-tasks map list"
#  2. The "TD category" sheet's column B ("Category Description" plus all
#     the per-row description text) is deleted outright, leaving only the
#     category names in column A. This also prunes the now-unused shared
#     strings, which is why so many <v> indices shift elsewhere in the
#     workbook - that's an automatic side effect, not a manual re-entry.
#  3. Various sheet selections move:
#       - "High level plan"            B3   -> B5
#       - "Activity-tasks map list"    A1:F25 -> C8
#       - "TD category"                B23  -> whole column B (B1:B1048576)
#  4. The active tab changes from "Technical debt tasks" (index 1) to
#     "TD category" (index 3) - i.e. "TD category" becomes the selected/
#     visible sheet when the workbook is reopened.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Activity-tasks map list" -> '"Activity-tasks map list"' ---
$wsMap = $wb.Worksheets.Item("Activity-tasks map list")
$wsMap.Name = '"Activity-tasks map list"'

# --- 3a. High level plan: move selection to B5 ---
$wsPlan = $wb.Worksheets.Item("High level plan")
$wsPlan.Range("B5").Select()

# --- 3b. Activity-tasks map list: move selection to C8 ---
$wsMap.Range("C8").Select()

# --- 2. TD category: drop column B (Category Description) entirely ---
$wsCat = $wb.Worksheets.Item("TD category")
$wsCat.Columns.Item(2).Delete()

# --- 3c. TD category: select column B (now empty) and make it the active sheet ---
$wsCat.Activate()
$wsCat.Range("B1:B1048576").Select()
